# "push before updating R" - append the latest working-hours log entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 20: Date / Hours / Notes, matching the formatting used by the
# existing log rows (date style on column A, wrapped text on column D).
$ws.Range("A20").NumberFormat = "d-mmm"
$ws.Range("A20").Value = Get-Date -Year 2024 -Month 5 -Day 22 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

$ws.Range("B20").Value = 5

$ws.Range("D20").Value = "Researching NMF as a method of grade estimation"
$ws.Range("D20").WrapText = $true

# Scroll the view down to the new row and leave the new note cell selected.
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("D20").Select()
